# Applies the "handles float input without breaking stuff" marksheet grading fix:
#  - Recomputes the Right/Wrong/Not-Attempt/Max summary (row 10) and the
#    Total row (row 12), fixing the marking-scheme cell (C11) which used to be
#    stored as text ("-1") and broke the total computation.
#  - Removes the unused third answer block (columns G:H) and the now-unused
#    portion of the second answer block (columns D:E, rows 19-40).
#  - Colors each student answer in column A (and, for the surviving rows,
#    column D) green/red depending on whether it matches the correct answer,
#    leaving it with the normal (black) style when the student left it blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellStyleFrom($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Row 9-12 summary block
# ---------------------------------------------------------------------------

# A10/A11/A12 ("No.", "Marking", "Total") get the mtitleStyle (same style as
# the row 9 header cells) instead of being left with the default style.
Set-CellStyleFrom "A9" "A10"
Set-CellStyleFrom "A9" "A11"
Set-CellStyleFrom "A9" "A12"

# Graded counts: 10 right, 4 wrong, 14 not attempted out of 28 questions.
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 28

# Marking scheme: +4 for right, -1 for wrong (now a real number, not text).
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Totals: 10*4 = 40, 4*-1 = -4, final score 36 out of max 112 (28*4).
$ws.Range("B12").Value = 40
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "36/112"

# ---------------------------------------------------------------------------
# Answer key blocks (rows 15-40)
# ---------------------------------------------------------------------------

# Drop the third "Student Ans / Correct Ans" block entirely (columns G:H).
$ws.Range("G15:H21").Clear()

# Drop the unused tail of the second block (columns D:E) for rows 19-40;
# rows 16-18 of that block stay and get the student's answers filled in.
$ws.Range("D19:E40").Clear()

Set-CellStyleFrom "B10" "D16"
$ws.Range("D16").Value = "Option A"
Set-CellStyleFrom "B10" "D17"
$ws.Range("D17").Value = "Option C"
Set-CellStyleFrom "B10" "D18"
$ws.Range("D18").Value = "Option D"

# Fill in column A (first block) with the student's answer, colored by
# correctness: correctStyle (green) when it matches column B, incorrectStyle
# (red) when it doesn't, and left as normalStyle/blank when not attempted.
Set-CellStyleFrom "B10" "A16"
$ws.Range("A16").Value = "Option A"
Set-CellStyleFrom "B10" "A17"
$ws.Range("A17").Value = "Option D"
Set-CellStyleFrom "B10" "A18"
$ws.Range("A18").Value = "Option B"
Set-CellStyleFrom "B10" "A22"
$ws.Range("A22").Value = "Option D"
Set-CellStyleFrom "C10" "A25"
$ws.Range("A25").Value = "Option D"
Set-CellStyleFrom "C10" "A28"
$ws.Range("A28").Value = "Option B"
Set-CellStyleFrom "B10" "A29"
$ws.Range("A29").Value = "Option D"
Set-CellStyleFrom "C10" "A31"
$ws.Range("A31").Value = "Option C"
Set-CellStyleFrom "B10" "A32"
$ws.Range("A32").Value = "Option C"
Set-CellStyleFrom "C10" "A37"
$ws.Range("A37").Value = "Option B"
Set-CellStyleFrom "B10" "A38"
$ws.Range("A38").Value = "Option A"
